$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.288.15"
$ws.Range("E2").Value = "  +2.03%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.098.11"
$ws.Range("E3").Value = "  +0.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.75%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.65"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("E6").Value = "  -0.64%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5281"
$ws.Range("E7").Value = "  +2.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4390"
$ws.Range("E8").Value = "  +0.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "55.08"
$ws.Range("E9").Value = "  +2.94%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09383"
$ws.Range("E10").Value = "  +2.13%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.175"
$ws.Range("E11").Value = "  +0.85%  "

# Row 12
$ws.Range("E12").Value = "  +0.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.565"
$ws.Range("E13").Value = "  +5.07%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.867"
$ws.Range("E14").Value = "  +1.52%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.985.39"
$ws.Range("E15").Value = "  -3.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.21"
$ws.Range("E16").Value = "  -1.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001159"
$ws.Range("E17").Value = "  +0.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -0.69%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.13"
$ws.Range("E19").Value = "  +0.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06728"
$ws.Range("E20").Value = "  +0.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.387"
$ws.Range("E21").Value = "  +3.06%  "

# Row 22
$ws.Range("E22").Value = "  -0.57%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.283.49"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.42"
$ws.Range("E24").Value = "  -1.70%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.989"
$ws.Range("E26").Value = "  +11.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.82"
$ws.Range("E27").Value = "  -0.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.50"
$ws.Range("E28").Value = "  +0.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.520"
$ws.Range("E29").Value = "  +1.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.81"
$ws.Range("E30").Value = "  +0.46%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.133"
$ws.Range("E31").Value = "  +0.43%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.688"
$ws.Range("E32").Value = "  +1.46%  "

# Row 33
$ws.Range("E33").Value = "  +0.34%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.255"
$ws.Range("E34").Value = "  +1.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.875"
$ws.Range("E35").Value = "  -2.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.08"
$ws.Range("E36").Value = "  -3.25%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02620"
$ws.Range("E37").Value = "  +1.82%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06767"
$ws.Range("E38").Value = "  +1.16%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.67"
$ws.Range("E39").Value = "  +1.95%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.352"
$ws.Range("E40").Value = "  +1.75%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6965"
$ws.Range("E41").Value = "  -0.21%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2216"
$ws.Range("E42").Value = "  +0.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6779"
$ws.Range("E43").Value = "  -0.32%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.30"
$ws.Range("E44").Value = "  +0.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.332"
$ws.Range("E45").Value = "  +0.71%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.311"
$ws.Range("E47").Value = "  +8.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.638"
$ws.Range("E48").Value = "  +0.43%  "

# Row 49
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.212"
$ws.Range("E49").Value = "  +6.46%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000340"
$ws.Range("E50").Value = "  -5.98%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07301"
$ws.Range("E51").Value = "  +3.52%  "
